# Applies the "Actualizada UD02 con inglés" edit:
#  1) Renames the two floating picture shapes:
#       image2.png -> image1.png
#       image1.png -> image3.png
#  2) Removes the whole "To deliver: entregar." bullet paragraph.
#  3) Updates the "To submit" gloss text to "enviar/entregar".

$d = $word.ActiveDocument

# --- 1) Rename floating picture shapes -------------------------------------
# Capture the original names up front so the two renames cannot collide
# (i.e. so the shape renamed to "image1.png" isn't then re-matched by the
# rule that looks for a shape still called "image1.png").
$originalShapeNames = @{}
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $originalShapeNames[$i] = $d.Shapes.Item($i).Name
}
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $oldName = $originalShapeNames[$i]
    if ($oldName -eq "image2.png") {
        $d.Shapes.Item($i).Name = "image1.png"
    } elseif ($oldName -eq "image1.png") {
        $d.Shapes.Item($i).Name = "image3.png"
    }
}

# --- 2) Delete the "To deliver: entregar." paragraph -----------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq "To deliver: entregar.") {
        $para.Range.Delete()
        break
    }
}

# --- 3) Update the "To submit" gloss ----------------------------------------
$d.Content.Find.Execute(
    " enviar (un fichero, un trabajo, etc.)", $true, $false, $false, $false,
    $false, $true, 1, $false,
    " enviar/entregar (un fichero, un trabajo, etc.)", 2)
